$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Binary Search analysis entry, mirroring the layout of row 3
$ws.Range("A4").Value = "Binary Search"
$ws.Range("D4").Value = "Type conversion"
$ws.Range("C4").Value = "We find mid element and compare it with target and based upon the comparision we reduce our search range."
$ws.Range("E4").Value = "O(log n)"
$ws.Range("F4").Value = "O(1)"
$ws.Range("B4").Value = "https://www.geeksforgeeks.org/binary-search/"

# Hyperlink for the resource URL in B4
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.geeksforgeeks.org/binary-search/")

# Styling: mirror row 3's layout (bold name in col A, wrapped notes in col C,
# plain text in col D, "Good" (green) highlight for the Big-O values)
$ws.Range("C4").Font.Bold = $false
$ws.Range("C4").WrapText = $true
$ws.Range("D4").Font.Bold = $false
$ws.Range("E4").Style = "Good"
$ws.Range("F4").Style = "Good"

$ws.Range("A6").Select()
